$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 119, shifting existing rows 119-200 down to 120-201.
$ws.Range("A119").EntireRow.Insert()

# Populate the new row 119 with its data.
$ws.Range("A119").Value = 9
$ws.Range("B119").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C119").Value = "Metropolitana"
$ws.Range("D119").Value = 44466
$ws.Range("E119").Value = 13
$ws.Range("F119").Value = 100112028
$ws.Range("G119").Value = "Sandia"
$ws.Range("H119").Value = "Sin especificar"
$ws.Range("I119").Value = "Primera"
$ws.Range("J119").Value = 97
$ws.Range("K119").Value = 1000
$ws.Range("L119").Value = 1000
$ws.Range("M119").Value = 1000
$ws.Range("N119").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O119").Value = "Perú"
$ws.Range("P119").Value = 1000
$ws.Range("Q119").Value = 1
$ws.Range("R119").Value = "Hortaliza"

# Match the date cell style used by the other rows in column D.
$ws.Range("D119").NumberFormat = $ws.Range("D120").NumberFormat
